$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")
if (-not $ws) { $ws = $wb.ActiveSheet }

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 24 de Octubre de 2020 a las 11:11"
$ws.Cells.Item(22, 2).Value = 385980
$ws.Cells.Item(22, 3).Value = 4070
$ws.Cells.Item(22, 4).Value = 309219
$ws.Cells.Item(22, 5).Value = 63556
$ws.Cells.Item(22, 7).Value = 128
$ws.Cells.Item(22, 8).Value = 13205
$ws.Cells.Item(23, 2).Value = 367819
$ws.Cells.Item(23, 3).Value = 2057
$ws.Cells.Item(23, 4).Value = 313112
$ws.Cells.Item(23, 5).Value = 47773
$ws.Cells.Item(23, 7).Value = 19
$ws.Cells.Item(23, 8).Value = 6934
$ws.Cells.Item(31, 1).Value = "Polonia"
$ws.Cells.Item(31, 2).Value = 241946
$ws.Cells.Item(31, 3).Value = 13628
$ws.Cells.Item(31, 4).Value = 109344
$ws.Cells.Item(31, 5).Value = 128251
$ws.Cells.Item(31, 7).Value = 179
$ws.Cells.Item(31, 8).Value = 4351
$ws.Cells.Item(32, 1).Value = "Chequia"
$ws.Cells.Item(32, 2).Value = 238323
$ws.Cells.Item(32, 4).Value = 91651
$ws.Cells.Item(32, 5).Value = 144701
$ws.Cells.Item(32, 8).Value = 1971
$ws.Cells.Item(66, 2).Value = 57965
$ws.Cells.Item(66, 3).Value = 14
$ws.Cells.Item(66, 5).Value = 105
$ws.Cells.Item(74, 1).Value = "Tunez"
$ws.Cells.Item(74, 2).Value = 48799
$ws.Cells.Item(74, 3).Value = 1585
$ws.Cells.Item(74, 4).Value = 5032
$ws.Cells.Item(74, 5).Value = 42948
$ws.Cells.Item(74, 7).Value = 35
$ws.Cells.Item(74, 8).Value = 819
$ws.Cells.Item(75, 1).Value = "Azerbaiyan"
$ws.Cells.Item(75, 2).Value = 48221
$ws.Cells.Item(75, 4).Value = 40831
$ws.Cells.Item(75, 5).Value = 6734
$ws.Cells.Item(75, 8).Value = 656
$ws.Cells.Item(76, 1).Value = "Kenia"
$ws.Cells.Item(76, 2).Value = 47843
$ws.Cells.Item(76, 4).Value = 33421
$ws.Cells.Item(76, 5).Value = 13538
$ws.Cells.Item(76, 8).Value = 884
$ws.Cells.Item(77, 1).Value = "Ghana"
$ws.Cells.Item(77, 2).Value = 47601
$ws.Cells.Item(77, 4).Value = 46824
$ws.Cells.Item(77, 5).Value = 463
$ws.Cells.Item(77, 8).Value = 314
$ws.Cells.Item(79, 1).Value = "Eslovaquia"
$ws.Cells.Item(79, 2).Value = 40801
$ws.Cells.Item(79, 3).Value = 2890
$ws.Cells.Item(79, 4).Value = 9920
$ws.Cells.Item(79, 5).Value = 30722
$ws.Cells.Item(79, 7).Value = 25
$ws.Cells.Item(79, 8).Value = 159
$ws.Cells.Item(80, 1).Value = "Afganistan"
$ws.Cells.Item(80, 2).Value = 40768
$ws.Cells.Item(80, 3).Value = 81
$ws.Cells.Item(80, 4).Value = 34023
$ws.Cells.Item(80, 5).Value = 5234
$ws.Cells.Item(80, 7).Value = 4
$ws.Cells.Item(80, 8).Value = 1511
$ws.Cells.Item(81, 1).Value = "Dinamarca"
$ws.Cells.Item(81, 2).Value = 38622
$ws.Cells.Item(81, 4).Value = 31295
$ws.Cells.Item(81, 5).Value = 6630
$ws.Cells.Item(81, 8).Value = 697
$ws.Cells.Item(82, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(82, 2).Value = 38493
$ws.Cells.Item(82, 4).Value = 26260
$ws.Cells.Item(82, 5).Value = 11168
$ws.Cells.Item(82, 8).Value = 1065
$ws.Cells.Item(83, 1).Value = "Serbia"
$ws.Cells.Item(83, 2).Value = 38115
$ws.Cells.Item(83, 4).Value = 31536
$ws.Cells.Item(83, 5).Value = 5793
$ws.Cells.Item(83, 8).Value = 786
$ws.Cells.Item(85, 1).Value = "Croacia"
$ws.Cells.Item(85, 2).Value = 33959
$ws.Cells.Item(85, 3).Value = 2242
$ws.Cells.Item(85, 4).Value = 23785
$ws.Cells.Item(85, 5).Value = 9745
$ws.Cells.Item(85, 7).Value = 16
$ws.Cells.Item(85, 8).Value = 429
$ws.Cells.Item(86, 1).Value = "El Salvador"
$ws.Cells.Item(86, 2).Value = 32585
$ws.Cells.Item(86, 3).Value = 164
$ws.Cells.Item(86, 4).Value = 28258
$ws.Cells.Item(86, 5).Value = 3383
$ws.Cells.Item(86, 7).Value = 4
$ws.Cells.Item(86, 8).Value = 944
$ws.Cells.Item(95, 2).Value = 19313
$ws.Cells.Item(95, 5).Value = 11440
$ws.Cells.Item(113, 2).Value = 9578
$ws.Cells.Item(113, 3).Value = 474
$ws.Cells.Item(113, 4).Value = 4015
$ws.Cells.Item(113, 5).Value = 5434
$ws.Cells.Item(113, 7).Value = 3
$ws.Cells.Item(113, 8).Value = 129
$ws.Cells.Item(142, 1).Value = "Estonia"
$ws.Cells.Item(142, 2).Value = 4351
$ws.Cells.Item(142, 3).Value = 52
$ws.Cells.Item(142, 4).Value = 3441
$ws.Cells.Item(142, 5).Value = 837
$ws.Cells.Item(142, 8).Value = 73
$ws.Cells.Item(143, 1).Value = "Islandia"
$ws.Cells.Item(143, 2).Value = 4308
$ws.Cells.Item(143, 4).Value = 3187
$ws.Cells.Item(143, 5).Value = 1110
$ws.Cells.Item(143, 8).Value = 11
$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 8).Value = 1
$ws.Cells.Item(217, 1).Value = "Islas Malvinas"
$ws.Cells.Item(217, 4).Value = 13
$ws.Cells.Item(217, 8).Value = 0
